# General enhancements and code clean-up.
# Adds a new "Project / Release" style column (E) to the Tardigrade version
# matrix, mirroring the existing B/C/D release columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column E. The order below matches the order the values were
# first typed by the original author (new row values, header last) so the
# shared-string table is appended to in the same sequence.
$ws.Range("E2").Value  = "9.0.0"
$ws.Range("E1").Value  = "11.0.0"
$ws.Range("E3").Value  = "5.0.0"
$ws.Range("E4").Value  = "3.3.0"
$ws.Range("E9").Value  = "1.3.0"
$ws.Range("E7").Value  = "9.3.0"
$ws.Range("E8").Value  = "8.2.0"
$ws.Range("E5").Value  = "2.0.0"
$ws.Range("E6").Value  = "6.0.0"
$ws.Range("E10").Value = "3.3.0"
$ws.Range("E11").Value = "1.1.0"

# Header row is bold, matching the rest of row 1.
$ws.Range("E1").Font.Bold = $true

# Match the (non-bold) formatting already used by the sibling C/D columns
# on the rows where they carry explicit formatting.
$ws.Range("E3").Font.Bold = $false
$ws.Range("E4").Font.Bold = $false
$ws.Range("E6").Font.Bold = $false
$ws.Range("E8").Font.Bold = $false
$ws.Range("E9").Font.Bold = $false
$ws.Range("E10").Font.Bold = $false
$ws.Range("E11").Font.Bold = $false

# The last cell touched/selected by the author.
$ws.Range("E8").Select()
